$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows at row 445, shifting existing data (old rows 445-480) down to 448-483
$ws.Rows("445:447").Insert()

$ws.Cells.Item(445,1).Value = 7
$ws.Cells.Item(445,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(445,3).Value = "Ñuble"
$ws.Cells.Item(445,4).Value = 44578
$ws.Cells.Item(445,5).Value = 16
$ws.Cells.Item(445,6).Value = "Fruta"
$ws.Cells.Item(445,7).Value = 100108
$ws.Cells.Item(445,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(445,9).Value = 100108006
$ws.Cells.Item(445,10).Value = "Plátano"
$ws.Cells.Item(445,11).Value = "Sin especificar"
$ws.Cells.Item(445,12).Value = "Pintón"
$ws.Cells.Item(445,13).Value = 240
$ws.Cells.Item(445,14).Value = 11000
$ws.Cells.Item(445,15).Value = 12000
$ws.Cells.Item(445,16).Value = 11500
$ws.Cells.Item(445,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(445,18).Value = "Ecuador"
$ws.Cells.Item(445,19).Value = 575
$ws.Cells.Item(445,20).Value = 20

$ws.Cells.Item(446,1).Value = 7
$ws.Cells.Item(446,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(446,3).Value = "Ñuble"
$ws.Cells.Item(446,4).Value = 44578
$ws.Cells.Item(446,5).Value = 16
$ws.Cells.Item(446,6).Value = "Fruta"
$ws.Cells.Item(446,7).Value = 100108
$ws.Cells.Item(446,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(446,9).Value = 100108006
$ws.Cells.Item(446,10).Value = "Plátano"
$ws.Cells.Item(446,11).Value = "Sin especificar"
$ws.Cells.Item(446,12).Value = "Primera"
$ws.Cells.Item(446,13).Value = 120
$ws.Cells.Item(446,14).Value = 15000
$ws.Cells.Item(446,15).Value = 15000
$ws.Cells.Item(446,16).Value = 15000
$ws.Cells.Item(446,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(446,18).Value = "Ecuador"
$ws.Cells.Item(446,19).Value = 750
$ws.Cells.Item(446,20).Value = 20

$ws.Cells.Item(447,1).Value = 7
$ws.Cells.Item(447,2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(447,3).Value = "Ñuble"
$ws.Cells.Item(447,4).Value = 44578
$ws.Cells.Item(447,5).Value = 16
$ws.Cells.Item(447,6).Value = "Fruta"
$ws.Cells.Item(447,7).Value = 100108
$ws.Cells.Item(447,8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(447,9).Value = 100108006
$ws.Cells.Item(447,10).Value = "Plátano"
$ws.Cells.Item(447,11).Value = "Sin especificar"
$ws.Cells.Item(447,12).Value = "Primera Pintón"
$ws.Cells.Item(447,13).Value = 120
$ws.Cells.Item(447,14).Value = 14000
$ws.Cells.Item(447,15).Value = 14000
$ws.Cells.Item(447,16).Value = 14000
$ws.Cells.Item(447,17).Value = "$/caja 20 kilos"
$ws.Cells.Item(447,18).Value = "Ecuador"
$ws.Cells.Item(447,19).Value = 700
$ws.Cells.Item(447,20).Value = 20

